# Insert a new weekly price record as the new first row of the Coliflor
# (Feria Lagunitas de Puerto Montt) data block. This pushes the existing
# rows 498:520 down to 499:521 and fills the freed row 498 with the new
# week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 498:520 down to 499:521, leaving row 498 empty for the new record.
$ws.Rows.Item(498).Insert()

# Populate the new row 498 with the new weekly observation.
$ws.Range("A498").Value = 4
$ws.Range("B498").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C498").Value = "Los Lagos"
$ws.Range("D498").Value = 45041
$ws.Range("E498").Value = 10
$ws.Range("F498").Value = 100112008
$ws.Range("G498").Value = "Coliflor"
$ws.Range("H498").Value = "Sin especificar"
$ws.Range("I498").Value = "Primera"
$ws.Range("J498").Value = 1200
$ws.Range("K498").Value = 1700
$ws.Range("L498").Value = 1700
$ws.Range("M498").Value = 1700
$ws.Range("N498").Value = "`$/unidad"
$ws.Range("O498").Value = "Región Metropolitana"
$ws.Range("P498").Value = 1700
$ws.Range("Q498").Value = 1
$ws.Range("R498").Value = "Hortaliza"
